$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns B:C down by one row for rows 2..36 -> 3..37, leaving column A
# untouched (the index column keeps its original row alignment). Work
# bottom-up so sources aren't clobbered before they're read.
for ($r = 36; $r -ge 2; $r--) {
    $dst = $r + 1
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2

    # Always clear the destination B:C first: Copy() from a blank source cell
    # does not clear an already-populated destination, so this emulates a
    # real shift-down for blank cells.
    $ws.Range("B$dst" + ":C$dst").ClearContents()

    # Re-populate using Copy() so the original type/shared-string-ness is
    # preserved (e.g. text that happens to look like a date, such as
    # "01/01/2020", stays literal text instead of becoming a date serial,
    # which is what a plain Value2 assignment would do). Only copy when
    # there's actually something there, otherwise Copy() leaves a stray
    # empty cell node behind instead of no node at all.
    if ($bVal -ne $null) {
        $ws.Range("B$r" + ":B$r").Copy($ws.Range("B$dst" + ":B$dst"))
    }
    if ($cVal -ne $null) {
        $ws.Range("C$r" + ":C$r").Copy($ws.Range("C$dst" + ":C$dst"))
    }
}

# New row 37 needs column A populated (the index sequence continues from 34
# to 35), using the same style as the rest of the index column: copy the
# format from A36, then overwrite with the new value.
$ws.Range("A36:A36").Copy($ws.Range("A37:A37"))
$ws.Cells.Item(37, 1).Value2 = 35

# Populate the newly-opened row 2 with the new "Derived Value" / "Value" pair.
$ws.Range("B2:C2").ClearContents()
$ws.Cells.Item(2, 2).Value2 = "Derived Value"
$ws.Cells.Item(2, 3).Value2 = "Value"

# The date that used to live in (old) row 2 / (new) row 3 gets a new value.
$ws.Cells.Item(3, 3).Value2 = "29/06/2020"
